$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = 6
    4  = 9
    6  = 10
    7  = 7
    8  = 9
    9  = 7
    10 = 7
    11 = 6
    13 = 9
    14 = 6
    15 = 6
    16 = 6
    17 = 6
    18 = 11
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
